# Update existing row 2 (ticker XBIT -> KB) and extend the transaction
# table with new rows (CIB, AVAL, a blank separator row, DMTK, GLBD and a
# trailing blank row), matching the "transaction builder" output shape:
# an index column (A) that always increments, plus a blank row every third
# record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Give the date column (B) a text format for the rows we are about to
# fill in, so "01/01/2020" is stored as literal text rather than being
# reinterpreted as a date serial number. ---
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B7").NumberFormat = "@"

function Set-TxRow($row, $date, $type, $ticker, $qty, $price, $cashflow, $cmlUnits, $cmlCost) {
    $ws.Cells.Item($row, 2).Value = $date          # B date
    $ws.Cells.Item($row, 3).Value = $type           # C type
    $ws.Cells.Item($row, 4).Value = $ticker         # D ticker
    $ws.Cells.Item($row, 5).Value = $qty             # E quantity
    $ws.Cells.Item($row, 6).Value = $price           # F price
    $ws.Cells.Item($row, 7).Value = 0                # G fees
    $ws.Cells.Item($row, 8).Value = 0                # H transact_val
    $ws.Cells.Item($row, 9).Value = 0                # I last_occurence
    $ws.Cells.Item($row, 10).Value = $cashflow       # J cashflow
    $ws.Cells.Item($row, 11).Value = 0               # K prev_units
    $ws.Cells.Item($row, 12).Value = $cmlUnits       # L cml_units
    $ws.Cells.Item($row, 13).Value = 0               # M prev_cost
    $ws.Cells.Item($row, 14).Value = $cmlCost        # N cml_cost
    $ws.Cells.Item($row, 15).Value = 0               # O cost_transact
    $ws.Cells.Item($row, 16).Value = 0               # P cost_unit
    $ws.Cells.Item($row, 17).Value = 0               # Q gain_loss
    $ws.Cells.Item($row, 18).Value = 0               # R yield
    $ws.Cells.Item($row, 19).Value = $price           # S avg_price
}

# Row 2 - retrained on the untrained year, ticker changed XBIT -> KB
Set-TxRow 2 "01/01/2020" "Buy" "KB"   120  41.36999893188477  -4964.399871826172  121.9512195121951  4964.399871826172

# Row 3 - new transaction (CIB)
Set-TxRow 3 "01/01/2020" "Buy" "CIB"  94   52.711181640625    -4954.85107421875   96.15384615384616  4954.85107421875

# Row 4 - new transaction (AVAL)
Set-TxRow 4 "01/01/2020" "Buy" "AVAL" 619  8.075420379638672  -4998.685214996338  625                4998.685214996338

# Row 5 - blank separator row produced by the (disorganized) transaction
# builder; only the running index column is populated.
$ws.Range("A4").Copy($ws.Range("A5"))
$ws.Cells.Item(5, 1).Value = 3

# Row 6 - new transaction (DMTK)
Set-TxRow 6 "01/01/2020" "Buy" "DMTK" 403  12.39999961853027  -4997.1998462677    416.6666666666667  4997.1998462677
$ws.Range("A4").Copy($ws.Range("A6"))
$ws.Cells.Item(6, 1).Value = 4

# Row 7 - new transaction (GLBD)
Set-TxRow 7 "01/01/2020" "Buy" "GLBD" 2941 1.700000047683716  -4999.700140237808  5000               4999.700140237808
$ws.Range("A4").Copy($ws.Range("A7"))
$ws.Cells.Item(7, 1).Value = 5

# Row 8 - trailing blank row, same pattern as row 5.
$ws.Range("A4").Copy($ws.Range("A8"))
$ws.Cells.Item(8, 1).Value = 6
